$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 620.1111
$ws.Range("I8").Value = 620.1111
$ws.Range("K8").Value = 1860.3333
$ws.Range("M8").Value = -1721.3333
$ws.Range("H80").Value = 2654.0667
$ws.Range("I80").Value = 1629.75
$ws.Range("J80").Value = 3026.5454
$ws.Range("K80").Value = 4889.25
$ws.Range("L80").Value = 9079.636200000001
$ws.Range("M80").Value = -3891.25
$ws.Range("N80").Value = -11075.6362
$ws.Range("H83").Value = 2654.0667
$ws.Range("I83").Value = 1629.75
$ws.Range("J83").Value = 3026.5454
$ws.Range("K83").Value = 14667.75
$ws.Range("L83").Value = 27238.9086
$ws.Range("M83").Value = -9675.75
$ws.Range("N83").Value = -37222.9086
$ws.Range("H86").Value = 7246.857
$ws.Range("I86").Value = 3332.6667
$ws.Range("J86").Value = 10182.5
$ws.Range("K86").Value = 3332.6667
$ws.Range("L86").Value = 10182.5
$ws.Range("M86").Value = -2209.6667
$ws.Range("N86").Value = -12428.5
$ws.Range("H88").Value = 2001136.2
$ws.Range("J88").Value = 2858116
$ws.Range("L88").Value = 2858116
$ws.Range("N88").Value = -2858928
$ws.Range("H89").Value = 7246.857
$ws.Range("I89").Value = 3332.6667
$ws.Range("J89").Value = 10182.5
$ws.Range("K89").Value = 16663.3335
$ws.Range("L89").Value = 50912.5
$ws.Range("M89").Value = -11047.3335
$ws.Range("N89").Value = -62144.5
$ws.Range("H91").Value = 2001136.2
$ws.Range("J91").Value = 2858116
$ws.Range("L91").Value = 2858116
$ws.Range("N91").Value = -2860924
$ws.Range("H94").Value = 3495
$ws.Range("I94").Value = 990
$ws.Range("J94").Value = 6000
$ws.Range("K94").Value = 990
$ws.Range("L94").Value = 6000
$ws.Range("M94").Value = -539
$ws.Range("N94").Value = -6902
$ws.Range("H112").Value = 475364.78
$ws.Range("J112").Value = 606694.5600000001
$ws.Range("L112").Value = 1820083.68
$ws.Range("N112").Value = -1822299.68
$ws.Range("H132").Value = 53753.95
$ws.Range("I132").Value = 53753.95
$ws.Range("K132").Value = 161261.85
$ws.Range("M132").Value = -158731.85
$ws.Range("H137").Value = 1749
$ws.Range("I137").Value = 1539.1111
$ws.Range("K137").Value = 4617.3333
$ws.Range("M137").Value = -2067.3333
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2250.3076
$ws.Range("I74").Value = 2161.2727
$ws.Range("K74").Value = 2161.2727
$ws.Range("M74").Value = -1287.2727
$ws.Range("H77").Value = 2250.3076
$ws.Range("I77").Value = 2161.2727
$ws.Range("K77").Value = 10806.3635
$ws.Range("M77").Value = -6438.363499999999
$ws.Range("H88").Value = 1778.4615
$ws.Range("J88").Value = 1645
$ws.Range("L88").Value = 1645
$ws.Range("N88").Value = -2457
$ws.Range("H91").Value = 1778.4615
$ws.Range("J91").Value = 1645
$ws.Range("L91").Value = 1645
$ws.Range("N91").Value = -4453
$ws.Range("H97").Value = 1809.7858
$ws.Range("I97").Value = 1803.4546
$ws.Range("K97").Value = 1803.4546
$ws.Range("M97").Value = -1307.4546
$ws.Range("H110").Value = 1531.625
$ws.Range("I110").Value = 1690.1666
$ws.Range("J110").Value = 1056
$ws.Range("K110").Value = 1690.1666
$ws.Range("L110").Value = 1056
$ws.Range("M110").Value = 354.8334
$ws.Range("N110").Value = -5146
$ws.Range("H132").Value = 55557320
$ws.Range("I132").Value = 1631.1666
$ws.Range("J132").Value = 166668700
$ws.Range("K132").Value = 4893.4998
$ws.Range("L132").Value = 500006100
$ws.Range("M132").Value = -2363.4998
$ws.Range("N132").Value = -500011160
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4713.5713
$ws.Range("J86").Value = 5199.6
$ws.Range("L86").Value = 5199.6
$ws.Range("N86").Value = -7445.6
$ws.Range("H89").Value = 4713.5713
$ws.Range("J89").Value = 5199.6
$ws.Range("L89").Value = 25998
$ws.Range("N89").Value = -37230
$ws.Range("H94").Value = 5758.04
$ws.Range("I94").Value = 1706.1052
$ws.Range("K94").Value = 1706.1052
$ws.Range("M94").Value = -1255.1052
$ws.Range("H134").Value = 41670692
$ws.Range("I134").Value = 19235102
$ws.Range("K134").Value = 57705306
$ws.Range("M134").Value = -57702771
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4547299.5
$ws.Range("I134").Value = 1887.6471
$ws.Range("J134").Value = 20001700
$ws.Range("K134").Value = 5662.9413
$ws.Range("L134").Value = 60005100
$ws.Range("M134").Value = -3127.9413
$ws.Range("N134").Value = -60010170
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 148.17857
$ws.Range("I38").Value = 85.85714
$ws.Range("J38").Value = 157.08163
$ws.Range("K38").Value = 257.57142
$ws.Range("L38").Value = 471.2448899999999
$ws.Range("M38").Value = 89.42858000000001
$ws.Range("N38").Value = -1165.24489
$ws.Range("H122").Value = 4036.3333
$ws.Range("J122").Value = 3277.5
$ws.Range("L122").Value = 29497.5
$ws.Range("N122").Value = -34397.5
$ws.Range("H131").Value = 492879.62
$ws.Range("I131").Value = 1033.2
$ws.Range("J131").Value = 675044.9399999999
$ws.Range("K131").Value = 3099.6
$ws.Range("L131").Value = 2025134.82
$ws.Range("M131").Value = 1940.4
$ws.Range("N131").Value = -2035214.82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1558.4
$ws.Range("I22").Value = 1558.4
$ws.Range("K22").Value = 1558.4
$ws.Range("M22").Value = -1263.4
$ws.Range("H27").Value = 1558.4
$ws.Range("I27").Value = 1558.4
$ws.Range("K27").Value = 1558.4
$ws.Range("M27").Value = -1451.4
$ws.Range("H40").Value = 3920
$ws.Range("I40").Value = 2785
$ws.Range("K40").Value = 2785
$ws.Range("M40").Value = -2649
$ws.Range("H42").Value = 11998.5
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 11998.5
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 11998.5
$ws.Range("M42").Value = ""
$ws.Range("N42").Value = -13124.5
$ws.Range("H46").Value = 2847.9
$ws.Range("H48").Value = 19832.334
$ws.Range("I48").Value = 19749.5
$ws.Range("J48").Value = 19998
$ws.Range("K48").Value = 19749.5
$ws.Range("L48").Value = 19998
$ws.Range("M48").Value = -19088.5
$ws.Range("N48").Value = -21320
$ws.Range("H49").Value = 11998.5
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 11998.5
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 11998.5
$ws.Range("M49").Value = ""
$ws.Range("N49").Value = -12292.5
$ws.Range("H68").Value = 2129.7
$ws.Range("J68").Value = 1651
$ws.Range("L68").Value = 1651
$ws.Range("N68").Value = -3149
$ws.Range("H71").Value = 2129.7
$ws.Range("J71").Value = 1651
$ws.Range("L71").Value = 8255
$ws.Range("N71").Value = -15743
$ws.Range("H82").Value = 2542.5715
$ws.Range("I82").Value = 2542.5715
$ws.Range("K82").Value = 2542.5715
$ws.Range("M82").Value = -2181.5715
$ws.Range("H85").Value = 2542.5715
$ws.Range("I85").Value = 2542.5715
$ws.Range("K85").Value = 2542.5715
$ws.Range("M85").Value = -1294.5715
$ws.Range("H93").Value = 83046.8
$ws.Range("I93").Value = 3616
$ws.Range("J93").Value = 136000.67
$ws.Range("K93").Value = 3616
$ws.Range("L93").Value = 136000.67
$ws.Range("M93").Value = -2368
$ws.Range("N93").Value = -138496.67
$ws.Range("H99").Value = 30000
$ws.Range("I99").Value = 30000
$ws.Range("K99").Value = 30000
$ws.Range("M99").Value = -27005
$ws.Range("H122").Value = 3477.4443
$ws.Range("I122").Value = 3196.5557
$ws.Range("J122").Value = 3758.3333
$ws.Range("K122").Value = 9589.667099999999
$ws.Range("L122").Value = 11274.9999
$ws.Range("M122").Value = -7139.667099999999
$ws.Range("N122").Value = -16174.9999
$ws.Range("H132").Value = 3088.5
$ws.Range("I132").Value = 3050.7778
$ws.Range("J132").Value = 3201.6667
$ws.Range("K132").Value = 9152.3334
$ws.Range("L132").Value = 9605.000100000001
$ws.Range("M132").Value = -6622.3334
$ws.Range("N132").Value = -14665.0001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9108
$ws.Range("J62").Value = 4000
$ws.Range("L62").Value = 4000
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 9108
$ws.Range("J65").Value = 4000
$ws.Range("L65").Value = 20000
$ws.Range("N65").Value = -26240
$ws.Range("H74").Value = 193591.75
$ws.Range("J74").Value = 193591.75
$ws.Range("L74").Value = 193591.75
$ws.Range("N74").Value = -195463.75
$ws.Range("H77").Value = 193591.75
$ws.Range("J77").Value = 193591.75
$ws.Range("L77").Value = 580775.25
$ws.Range("N77").Value = -590135.25
